# Clean up / fix emission factors on the "emission_factors" sheet.
#
# 1) Rows 81-131 had an H column formula ("=I##*1.15" / "=I##*1.12") that
#    referenced the now-redundant helper column I. Freeze those formulas to
#    their last computed value, then delete column I outright (its data was
#    only ever a scratch/helper column for the H formulas and a couple of
#    stray notes).
# 2) A handful of placeholder 0 values (rows 146-151) get their real,
#    previously-missing emission factor values filled in.
# 3) Turn on AutoFilter for the table so it's easier to work with.
# 4) Leave the cursor/selection near the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("emission_factors")
$ws.Activate()

# --- 1) Freeze H81:H131 formulas to static values, then drop column I ----
for ($r = 81; $r -le 131; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = $cell.Value2
}

# Any leftover comment anchored in column I goes away with the column.
foreach ($cmt in @("I87")) {
    $c = $ws.Range($cmt).Comment
    if ($c -ne $null) {
        $c.Delete()
    }
}

$ws.Columns.Item(9).Delete()

# --- 2) Fill in the real values that used to be placeholder zeroes -------
$ws.Range("H146").Value = 1.359
$ws.Range("H147").Value = 1.237
$ws.Range("H148").Value = 1.509
$ws.Range("H149").Value = 0.817
$ws.Range("H150").Value = 0.971
$ws.Range("H151").Value = 0.629

# --- 3) Turn on AutoFilter over the real table extent ---------------------
$ws.Range("A1:H355").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=emission_factors!`$A`$1:`$H`$355")
$filterName.Visible = $false

# --- 4) Park the selection near the top-left of the sheet -----------------
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H1").Select()
